$d = $word.ActiveDocument
$endPos = $d.Content.End
$r = $d.Range($endPos, $endPos)
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p/><w:p><w:r><w:t>9:</w:t></w:r></w:p><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Igen recon10 = fixed og recon40 = movin</w:t></w:r><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">g. </w:t></w:r><w:r><w:t xml:space="preserve">Trækker gennemsnit fra, clamper, </w:t></w:r><w:r><w:t>o</w:t></w:r><w:r><w:t xml:space="preserve">g normaliserer til -1000 og 1000. Har ikke lagt datafilerne op men preprocess scriptet, så jeg kan se det bedre. </w:t></w:r><w:r><w:t xml:space="preserve">Nu </w:t></w:r><w:r><w:t xml:space="preserve">har jeg valgt </w:t></w:r><w:r><w:t>(NumberOfSpatialSamples 100000)</w:t></w:r><w:r><w:t xml:space="preserve"> i stedet for </w:t></w:r><w:r><w:t xml:space="preserve">(NumberOfSpatialSamples </w:t></w:r><w:r><w:t>5</w:t></w:r><w:r><w:t>00000)</w:t></w:r><w:r><w:t xml:space="preserve">. </w:t></w:r><w:r><w:t xml:space="preserve">Se på screenshottet af sammenligning mellem result (venstre) og fixed (højre). </w:t></w:r><w:r><w:t>Det er egentlig rigtig nok, men den skubber billedet lidt.. Mærkeligt…</w:t></w:r></w:p><w:p/><w:p><w:r><w:t>10:</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Prøvede med præcis samme som ovenfor for at teste, om de sorte prikker </w:t></w:r><w:r><w:t xml:space="preserve">(meget negative værdier) kommer pga. </w:t></w:r><w:r><w:t xml:space="preserve">jeg ændrede </w:t></w:r><w:r><w:t>NumberOfSpatialSamples</w:t></w:r><w:r><w:t xml:space="preserve"> til 100000 i stedet for 500000. Så kører præcis samme som ovenfor, men med 500000 i stedet for </w:t></w:r><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">100000 og det giver præcis det samme. Det er nærmest helt ens. Jeg fik de negative tal i workprogress 7, så tænker, at det måske skyldes </w:t></w:r><w:r><w:t>primært, at jeg ikke normaliserer.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r.InsertXML($xml)
